# Fill in the previously blank table row (row 11) with the new
# "Log - erro de login" requirement (functional + technical description),
# matching the new table row added to the "Lista de verificação" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "Log - erro de login"
$ws.Range("C11").Value = "Registra cada tentativa de login falha, incluindo data e hora do evento, arquitetura da máquina, sistema operacional, e a mensagem de erro."
$ws.Range("D11").Value = "Verifica as credenciais fornecidas com as armazenadas no banco de dados, se não forem as mesmas será criado um aquivo txt   através da aplicação java. Será criado um arquivo txt por dia, o nome será o ano, mês e dia, do dia criado."

# Highlight the new requirement name in bold, same as the other
# requirement-name cells in column B.
$ws.Range("B11").Font.Bold = $true

# Let the newly-typed, wrapped text re-flow the row heights (and the
# rows above it whose wrapped text reflows once the font table changes).
$ws.Rows.Item(5).RowHeight = 96.6
$ws.Rows.Item(6).RowHeight = 151.8
$ws.Rows.Item(11).RowHeight = 82.8

# Move the active selection the way the author left it after finishing
# the edit.
$ws.Range("H11").Select()
